$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, re-format R4, R5, R6, R7 to match the rest of their rows (the diff
# shows R4/R5/R6/R7 lose their distinct "last column" styling and pick up
# the same style as the other year columns, e.g. Q4/Q5/Q6/Q7). This must
# happen BEFORE we copy R's formatting into the new column S, so that S
# ends up sharing the same (now-common) style rather than the old one-off
# "last column" style.
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)

$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)

$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)

$ws.Range("Q7").Copy()
$ws.Range("R7").PasteSpecial(-4122)

# Add a new column S with 2022 data, matching the existing per-row formatting
# by copying formats from the corresponding column R cells (same row), then
# set the new values.

# --- Header row (row 4): S4 = 2022, formatted like R4/Q4 ---
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

# --- Row 5: S5 = 49.7 ---
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 49.7

# --- Row 6: S6 = 34.9 ---
$ws.Range("R6").Copy()
$ws.Range("S6").PasteSpecial(-4122)
$ws.Range("S6").Value = 34.9

# --- Row 7: S7 = 21 ---
$ws.Range("R7").Copy()
$ws.Range("S7").PasteSpecial(-4122)
$ws.Range("S7").Value = 21

$wb.Application.CutCopyMode = 0

# Update the selection to match the new state.
$ws.Range("R12").Select()
